# Reorder the player roster table (rows 2-18) on Sheet1.
# Each row's Position and Team values travel together with the Player
# so that the underlying player/position/team association is preserved;
# only the row order (as stored) changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Mikal Bridges",      "SG,SF,PF", "New York Knicks"),
    @("Amen Thompson",      "SG,SF",    "Houston Rockets"),
    @("Santi Aldama",       "PF,C",     "Memphis Grizzlies"),
    @("Brook Lopez",        "C",        "Milwaukee Bucks"),
    @("Nikola Vucevic",     "PF,C",     "Chicago Bulls"),
    @("Evan Mobley",        "PF,C",     "Cleveland Cavaliers"),
    @("Bennedict Mathurin", "SG,SF",    "Indiana Pacers"),
    @("Josh Giddey",        "PG,SG,SF", "Chicago Bulls"),
    @("DeMar DeRozan",      "SF,PF",    "Sacramento Kings"),
    @("Kelly Oubre Jr.",    "SG,SF",    "Philadelphia 76ers"),
    @("Scottie Barnes",     "SG,SF,PF", "Toronto Raptors"),
    @("Luka Doncic",        "PG,SG",    "Dallas Mavericks"),
    @("De'Aaron Fox",       "PG",       "Sacramento Kings"),
    @("Tyler Herro",        "PG,SG",    "Miami Heat"),
    @("Caris LeVert",       "SG,SF",    "Cleveland Cavaliers"),
    @("Miles Bridges",      "SF,PF",    "Charlotte Hornets"),
    @("Ja Morant",          "PG",       "Memphis Grizzlies")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
